# Update column G ("K" - strikeouts) values on Sheet1 to reflect the
# regenerated save_data (K instead of Strike#), per rows 2-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 4
    3  = 1
    4  = 1
    5  = 4
    6  = 5
    7  = 0
    8  = 0
    9  = 0
    10 = 4
    11 = 1
    12 = 1
    13 = 0
    14 = 5
    15 = 2
    16 = 3
    17 = 4
    18 = 4
    19 = 4
    20 = 4
    21 = 4
    22 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
